# Power vs. Throughput - apply "Updated power vs throughput" commit
$d = $word.ActiveDocument

# Unicode helpers for curly quotes used in the new text
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”

# ---------------------------------------------------------------------
# Paragraph: "This familiar regulation makes sense ..."
#   "makes sense" -> "makes solid sense"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "regulation makes sense in a world of analog communications.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "regulation makes solid sense in a world of analog communications.",
    2) | Out-Null

# ---------------------------------------------------------------------
# Paragraph: "More operators can be accommodated ..."
#   "Power affects bandwidth occupancy." -> "Power used affects the bandwidth occupied."
#   "have more people communicating, and that is a very worthy goal."
#     -> "have more people communicating per unit time, and that is a very
#         worthy goal, and that is why this regulation exists. "
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Power affects bandwidth occupancy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Power used affects the bandwidth occupied.",
    2) | Out-Null

$d.Content.Find.Execute(
    "have more people communicating, and that is a very worthy goal.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "have more people communicating per unit time, and that is a very worthy goal, and that is why this regulation exists. ",
    2) | Out-Null

# ---------------------------------------------------------------------
# Paragraph: "If we take the goal of more people being able to communicate..."
#   "If we take the goal of more people being able to communicate, and we"
#     -> 'If we take on the goal of "more people being able to communicate
#         per unit time", and we'
#   trailing space added after "conserve?"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "If we take the goal of more people being able to communicate, and we",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If we take on the goal of " + $ldq + "more people being able to communicate per unit time" + $rdq + ", and we",
    2) | Out-Null

$d.Content.Find.Execute(
    "power that we should require operators to conserve?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "power that we should require operators to conserve? ",
    2) | Out-Null

# ---------------------------------------------------------------------
# Insert four new paragraphs right before the bookmark paragraph
# (the one holding the _GoBack bookmark, originally the last, empty,
# paragraph in the body).
# ---------------------------------------------------------------------
$newParaTexts = @(
    "It may not be. In order to achieve the goal of more people communicating in a digital channel, we must confront the idea of throughput. ",
    "Optimizing throughput in the analog realm means minimizing power. Optimizing throughput in the digital realm means that we have to consider both power and coding gain. ",
    "Digital signals generally receive one or both of the following types of coding. A signal is sampled, and converted into a series of discrete numbers that represent the signal. Once we have that set of numbers, we can remove unnecessary redundancy. This is, essentially, compression. After that, we add the right type of redundancy in order to make the signal resilient to all the things it will encounter as it" + [char]0x2019 + "s sent over the air. This is forward error correction coding. It" + [char]0x2019 + "s like armor. The end result of all this coding is very effective gain. Our signal acts like it" + [char]0x2019 + "s much more powerful. We can receive it and reconstruct a perfect or almost-perfect version of the transmitted signal. This signal can be sent at a lower power than an equivalent analog signal because we can use math to help fix errors caused by noise or interference. We don" + [char]0x2019 + "t have to bull our way through the static with a loud analog signal. We can correct errors. "
)

$bmParaIndex = $d.Paragraphs.Count
foreach ($t in $newParaTexts) {
    $bmPara = $d.Paragraphs.Item($bmParaIndex)
    $insPoint = $d.Range($bmPara.Range.Start, $bmPara.Range.Start)
    $insPoint.InsertParagraphBefore()
    $bmParaIndex = $bmParaIndex + 1
    $newPara = $d.Paragraphs.Item($bmParaIndex - 1)
    $newPara.Range.InsertBefore($t)
}

# The bookmark paragraph itself gains a leading sentence before the bookmark:
# "An additional advantage is that with some other math tricks, we can
#  control the occupied bandwidth. "
$bmPara = $d.Paragraphs.Item($bmParaIndex)
$bmInsertPoint = $d.Range($bmPara.Range.Start, $bmPara.Range.Start)
$bmInsertPoint.InsertBefore("An additional advantage is that with some other math tricks, we can control the occupied bandwidth. ")

# ---------------------------------------------------------------------
# New paragraph after the bookmark paragraph, plus a trailing blank
# paragraph before the section break.
# ---------------------------------------------------------------------
$bmPara = $d.Paragraphs.Item($bmParaIndex)
$bmPara.Range.InsertParagraphAfter()
$decidingParaIndex = $bmParaIndex + 1
$decidingPara = $d.Paragraphs.Item($decidingParaIndex)
$decidingPara.Range.InsertBefore("Deciding how and when to twiddle two knobs (power and coding gain) instead of one knob (power) means the complexity has increased. ")
$decidingPara.Range.InsertParagraphAfter()
